# Fix units in stages tables.
# The "lm-hr" / "kW-hr" unit labels in the Lamp rows (2, 4, 6, 8) of the
# "EEU data" sheet are replaced with the macro names ("\lmhr" / "\kWhr")
# that actually expand into those strings (with the missing dot) in the
# ReboundPaper2022 / ReboundTools2022 repositories. Threaded comments are
# left on E2/F2 explaining the macros, and the selection is moved to F2.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = "\lmhr"
$ws.Range("F2").Value = "\kWhr"
$ws.Range("E4").Value = "\lmhr"
$ws.Range("F4").Value = "\kWhr"
$ws.Range("E6").Value = "\lmhr"
$ws.Range("F6").Value = "\kWhr"
$ws.Range("E8").Value = "\lmhr"
$ws.Range("F8").Value = "\kWhr"

$ws.Range("E2").AddCommentThreaded("\lmhr is a macro in the ReboundPaper2022 repository hat expands into the right thing (lm-hr but with a dot)") | Out-Null
$ws.Range("F2").AddCommentThreaded("\kWhr is a macro in the ReboundTools2022 repository that expands into the right thing (kW-hr but with a dot)") | Out-Null

$ws.Range("F2").Select() | Out-Null
